{"js": "// Replace each division expression in the table with its updated value.\n// Old -> New values mirror the canonical diff (each <w:t> run's text is\n// swapped in place); matchCase keeps the search exact since the table also\n// holds the header date text.\nconst replacements = [\n  [\"976\u00f77=\", \"510\u00f77=\"],\n  [\"408\u00f73=\", \"134\u00f78=\"],\n  [\"217\u00f74=\", \"930\u00f75=\"],\n  [\"722\u00f75=\", \"205\u00f73=\"],\n  [\"628\u00f76=\", \"109\u00f74=\"],\n  [\"968\u00f77=\", \"753\u00f76=\"],\n  [\"674\u00f74=\", \"848\u00f73=\"],\n  [\"194\u00f79=\", \"340\u00f77=\"],\n  [\"651\u00f79=\", \"369\u00f76=\"],\n  [\"177\u00f75=\", \"950\u00f78=\"],\n  [\"724\u00f75=\", \"793\u00f77=\"],\n  [\"864\u00f79=\", \"538\u00f73=\"],\n  [\"838\u00f77=\", \"928\u00f76=\"],\n  [\"539\u00f73=\", \"663\u00f76=\"],\n  [\"854\u00f79=\", \"415\u00f72=\"],\n  [\"122\u00f73=\", \"712\u00f74=\"],\n  [\"603\u00f79=\", \"744\u00f78=\"],\n  [\"464\u00f73=\", \"536\u00f79=\"],\n  [\"380\u00f72=\", \"595\u00f72=\"],\n  [\"818\u00f79=\", \"141\u00f76=\"],\n  [\"100\u00f77=\", \"711\u00f77=\"],\n  [\"888\u00f75=\", \"879\u00f77=\"],\n  [\"176\u00f79=\", \"855\u00f72=\"],\n  [\"234\u00f72=\", \"245\u00f77=\"],\n  [\"109\u00f77=\", \"742\u00f74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const found = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each division expression in the table with its updated value.\n# Old -> New values are applied via Find/Replace across the whole document body,\n# matching the canonical diff (each <w:t> run's text is swapped in place).\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"976\u00f77=\", \"510\u00f77=\"),\n    @(\"408\u00f73=\", \"134\u00f78=\"),\n    @(\"217\u00f74=\", \"930\u00f75=\"),\n    @(\"722\u00f75=\", \"205\u00f73=\"),\n    @(\"628\u00f76=\", \"109\u00f74=\"),\n    @(\"968\u00f77=\", \"753\u00f76=\"),\n    @(\"674\u00f74=\", \"848\u00f73=\"),\n    @(\"194\u00f79=\", \"340\u00f77=\"),\n    @(\"651\u00f79=\", \"369\u00f76=\"),\n    @(\"177\u00f75=\", \"950\u00f78=\"),\n    @(\"724\u00f75=\", \"793\u00f77=\"),\n    @(\"864\u00f79=\", \"538\u00f73=\"),\n    @(\"838\u00f77=\", \"928\u00f76=\"),\n    @(\"539\u00f73=\", \"663\u00f76=\"),\n    @(\"854\u00f79=\", \"415\u00f72=\"),\n    @(\"122\u00f73=\", \"712\u00f74=\"),\n    @(\"603\u00f79=\", \"744\u00f78=\"),\n    @(\"464\u00f73=\", \"536\u00f79=\"),\n    @(\"380\u00f72=\", \"595\u00f72=\"),\n    @(\"818\u00f79=\", \"141\u00f76=\"),\n    @(\"100\u00f77=\", \"711\u00f77=\"),\n    @(\"888\u00f75=\", \"879\u00f77=\"),\n    @(\"176\u00f79=\", \"855\u00f72=\"),\n    @(\"234\u00f72=\", \"245\u00f77=\"),\n    @(\"109\u00f77=\", \"742\u00f74=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 1, $false, $pair[1], 2)\n}\n"}
